$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
  2 = @{"B"="10.52666939636015"; "C"="9.054858149168012"; "E"="16.26962362656793"; "F"="16.86991607391233"; "G"="3.687180368146247"; "I"="32.63602809777574"; "K"="11.89819788035099"}
  3 = @{"B"="10.36722544794065"; "C"="8.736599914076832"; "E"="15.35947500710867"; "F"="15.89584955866808"; "G"="3.690685420601557"; "I"="32.119942038948"; "K"="11.75815151589214"}
  4 = @{"B"="10.27286681707926"; "C"="8.537844792923561"; "E"="14.7777930532572"; "F"="15.26997757108491"; "G"="3.69294296036939"; "I"="31.8010068668614"; "K"="11.67565650702839"}
  5 = @{"B"="10.23536489433787"; "C"="8.456143148135002"; "E"="14.5352840896444"; "F"="15.00819731993403"; "G"="3.693889550902492"; "I"="31.67062624487979"; "K"="11.64295873680934"}
  6 = @{"B"="10.22919680489135"; "C"="8.442538043682463"; "E"="14.49469446763265"; "F"="14.96433081551593"; "G"="3.69404834298294"; "I"="31.64895484617255"; "K"="11.63758599494237"}
  7 = @{"B"="10.27235712949926"; "C"="8.536745615573668"; "E"="14.77454423173287"; "F"="15.26647399323137"; "G"="3.692955618468436"; "I"="31.79925003694136"; "K"="11.6752117575025"}
  8 = @{"B"="10.47099471666843"; "C"="8.945893208708929"; "E"="15.96069639732649"; "F"="16.53996406344768"; "G"="3.688367099185589"; "I"="32.45857116020861"; "K"="11.84920996145398"}
  9 = @{"B"="10.88590085727118"; "C"="9.716410415222509"; "E"="18.10541171851309"; "F"="19.00274580682531"; "G"="3.680200238989138"; "I"="33.73063383505158"; "K"="12.21624412587892"}
  10 = @{"B"="11.20249500382399"; "C"="10.25682024703986"; "E"="19.72590634603953"; "F"="20.67494806633232"; "G"="3.674699361316546"; "I"="34.64655190349175"; "K"="12.49899878719998"}
  11 = @{"B"="11.34828729605643"; "C"="10.49597935423284"; "E"="20.42214917594453"; "F"="21.3917225636224"; "G"="3.67230369729545"; "I"="35.05796437956845"; "K"="12.62989342896227"}
  12 = @{"B"="11.40368532875158"; "C"="10.58550414858688"; "E"="20.67995164483167"; "F"="21.65686569030329"; "G"="3.671411743424828"; "I"="35.21290454612928"; "K"="12.67973701746218"}
  13 = @{"B"="11.39174695560266"; "C"="10.56627078264433"; "E"="20.62468886293896"; "F"="21.60004134736742"; "G"="3.671603165998047"; "I"="35.17957486983175"; "K"="12.66899079063549"}
  14 = @{"B"="11.35284149180559"; "C"="10.50336586911309"; "E"="20.44347576966054"; "F"="21.4136618050453"; "G"="3.672230011087673"; "I"="35.070729044082"; "K"="12.63398882944611"}
  15 = @{"B"="11.32903356517829"; "C"="10.46469721339704"; "E"="20.33171679116245"; "F"="21.29868154950795"; "G"="3.67261595211718"; "I"="35.00394390533246"; "K"="12.61258367824585"}
  16 = @{"B"="11.19299718590051"; "C"="10.24104927533424"; "E"="19.67958472674569"; "F"="20.62722412089977"; "G"="3.674858061263754"; "I"="34.61955108617422"; "K"="12.49048591551105"}
  17 = @{"B"="11.10995028277924"; "C"="10.10207890583352"; "E"="19.26906567543916"; "F"="20.20408069617459"; "G"="3.676260773865155"; "I"="34.38232408705281"; "K"="12.41612850864091"}
  18 = @{"B"="11.06235450788991"; "C"="10.02152256083806"; "E"="19.02909147484874"; "F"="19.95656407809808"; "G"="3.677077628133659"; "I"="34.24538955202871"; "K"="12.3735765426531"}
  19 = @{"B"="11.04627070929097"; "C"="9.994142932409382"; "E"="18.9471772249569"; "F"="19.87204792380562"; "G"="3.677355930678804"; "I"="34.19894515877935"; "K"="12.35920786609869"}
  20 = @{"B"="11.11877358734593"; "C"="10.11693775490017"; "E"="19.31316479497694"; "F"="20.2495528364879"; "G"="3.67611041328568"; "I"="34.40762852387608"; "K"="12.42402193966674"}
  21 = @{"B"="11.36426433333394"; "C"="10.52187136722456"; "E"="20.49686091879133"; "F"="21.46857628470567"; "G"="3.672045478996158"; "I"="35.10272361479862"; "K"="12.64426263934269"}
  22 = @{"B"="11.52578324948083"; "C"="10.78042920001486"; "E"="21.23639546055286"; "F"="22.22866616901555"; "G"="3.669477541289486"; "I"="35.55199460043686"; "K"="12.7897924093039"}
  23 = @{"B"="11.43950013171811"; "C"="10.64301306587373"; "E"="20.84479841004908"; "F"="21.82633154475864"; "G"="3.670840016370481"; "I"="35.31270061123987"; "K"="12.71199122965749"}
  24 = @{"B"="11.11478410230026"; "C"="10.11022212018468"; "E"="19.29323992440774"; "F"="20.22900810905294"; "G"="3.676178358851217"; "I"="34.39619009224631"; "K"="12.42045270182471"}
  25 = @{"B"="10.77132237292633"; "C"="9.512064450417725"; "E"="17.53979093795894"; "F"="18.34778573295697"; "G"="3.682321370716225"; "I"="33.38936394186408"; "K"="12.11445622438158"}
}

foreach ($row in $data.Keys) {
  foreach ($col in $data[$row].Keys) {
    $ws.Range("$col$row").Value = [double]$data[$row][$col]
  }
}
